# The sheet was originally exported by pandas with an extra leading
# "index" column (A: 0,1,2,...) and an extra leading metadata row
# (row 1: numeric column positions 0,2) baked in. This commit removes
# both of those pandas artifacts so the real header/data starts at A1.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(1).EntireRow.Delete()
$ws.Columns.Item(1).EntireColumn.Delete()
